$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Flights")

$ws.Range("A8").Value = "Milan"
$ws.Range("A9").Value = "Madrid"
$ws.Range("A10").Value = "Barcelona"
$ws.Range("A11").Value = "Sydney"
$ws.Range("A12").Value = "Newcastle"

$ws.Activate()
$ws.Range("A13").Select()
